$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

# Row 8: B8 becomes a numeric value (was the "Selenium" text); row height shrinks slightly
$ws.Range("B8").Value = 12
$ws.Rows.Item(8).RowHeight = 13.8

# Row 9 (new row): A9 holds the renamed string (was "Selenium", now "Search Query 2"),
# B9 holds a numeric value
$ws.Range("A9").Value = "Search Query 2"
$ws.Range("B9").Value = 12.2

# Update selection to match the target state
$ws.Range("B10").Select()
